# Daily attendance processing - 2026-01-12 09:14:20
#
# The "Recorded By" column (G) lists the recorder(s) for each session as a
# comma-separated string. Rows where the session was completed automatically
# show "<user email>, System"; normalize these to lead with "System" instead,
# i.e. "System, <user email>".
#
# Scan every used row in column G and swap the two comma-separated tokens
# whenever the value is exactly "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
        $changed = $changed + 1
    }
}

Write-Output ("Updated 'Recorded By' cells: " + $changed)
